# Apply the JD-sample edit described by the commit "Data regarding skill added".
#
# Three content-level changes plus one style-definition tweak:
#   1. "Company:" line - the " " / "TechNova" / " Solutions " runs (with the
#      spell-check proofErr wrapper around "TechNova") collapse into a single
#      run " TechNova Solutions " once the same text is found & replaced.
#   2. The "Technical:" bullet - drop the spell-check proofErr wrapper that
#      surrounds the bold "PyTorch" run, without touching the run itself.
#   3. The "Experience:" bullet - the trailing "." run (sibling of the bold
#      "Kubernetes" run) becomes ", nlp." - i.e. the "nlp" skill is appended.
#   4. styles.xml - "Default Paragraph Font" gets flagged semi-hidden.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Company: <space>TechNova<space>Solutions<space>
#    Re-typing the identical visible text over the three runs (the middle
#    one sandwiched between <w:proofErr w:type="spellStart"/> / spellEnd)
#    makes Word re-emit it as one plain run and drops the now-stale
#    proofErr markers around "TechNova".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" TechNova Solutions ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " TechNova Solutions ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Technical: ... PyTorch ...
#    "PyTorch" is its own (bold) run with a proofErr spellStart/spellEnd
#    pair immediately outside of it. Widen a range one character to each
#    side (so both proofErr markers fall strictly inside it), capture the
#    plain text, delete the whole span (which drops the enclosed proofErr
#    markers) and retype it, then restore the bold run/formatting on just
#    the "PyTorch" word so the rest of the sentence is unaffected.
# ---------------------------------------------------------------------
$pyRng = $d.Content
$pyRng.Find.Execute("PyTorch") | Out-Null
$pyRng.MoveStart(1, -1)
$pyRng.MoveEnd(1, 1)
$pyText = $pyRng.Text
$pyRng.Delete()
$pyRng.InsertAfter($pyText)

$pyWord = $d.Content
$pyWord.Find.Execute("PyTorch") | Out-Null
$pyWord.Bold = 1

# ---------------------------------------------------------------------
# 3) Experience: ... Kubernetes.
#    Only the trailing "." run (a sibling run after the bold "Kubernetes"
#    run) changes, becoming ", nlp." - the bold "Kubernetes" run itself
#    must stay untouched.
# ---------------------------------------------------------------------
$kubeRng = $d.Content
$kubeRng.Find.Execute("Kubernetes.") | Out-Null
$dotRng = $d.Range($kubeRng.End - 1, $kubeRng.End)
$dotRng.Text = ", nlp."

# ---------------------------------------------------------------------
# 4) styles.xml - mark the built-in "Default Paragraph Font" character
#    style semi-hidden (as happens when the style gets hidden from the
#    style gallery but left available via "unhide when used").
# ---------------------------------------------------------------------
$dpf = $d.Styles("Default Paragraph Font")
$dpf.UnhideWhenUsed = $true
$dpf.Visibility = $false
